$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.518.05"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.996.52"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.13"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.96"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.995.14"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.89"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.11%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.23"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.64%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.491.57"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.04"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.510.07"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.990.78"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "453.75"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.98"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.35"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.97"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.23"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.49"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.03"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.01%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.42%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.16"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.49"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.108"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0832"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.93%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.79"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "9.28"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.99%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.23"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.06"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.47%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.123"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +8.29%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.88"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "394.99"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -7.06%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0353"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Arweave"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "39.03"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.81%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.717.93"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.41"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.62%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.108"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.30%  "
